# Fix mislabeled property_category values:
#  - Sheet "建物" (Building): rows 2-6, column I (property_category) were "land" -> should be "building"
#  - Sheet "汽車" (Car): row 2, column H (property_category) was "land" -> should be "car"

$wb = $excel.ActiveWorkbook

$wsBuilding = $wb.Worksheets.Item("建物")
for ($r = 2; $r -le 6; $r++) {
    $wsBuilding.Cells.Item($r, 9).Value = "building"
}

$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Cells.Item(2, 8).Value = "car"
